$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "NewMotor" / "NewPE" columns (E, F) ---------------------------------
# Header row (bold, reuse existing header style)
$ws.Range("E1").Value = "NewMotor"
$ws.Range("F1").Value = "NewPE"
$ws.Range("E1:F1").Font.Bold = $true

# Data rows 2-11
$ePE = @{
    2  = @("LiftLeftA",  "-")
    3  = @("DriveLeftB", "D")
    4  = @("DriveRightB","C")
    5  = @("LiftRightA", "-")
    6  = @("LiftLeftB",  "-")
    7  = @("LiftRightB", "-")
    8  = @("LiftLeftC",  "-")
    9  = @("DriveLeftA", "B")
    10 = @("DriveRightA","A")
    11 = @("LiftRightC", "-")
}

foreach ($r in 2..11) {
    $vals = $ePE[$r]
    $ws.Cells.Item($r, 5).Value = $vals[0]
    $ws.Cells.Item($r, 6).Value = $vals[1]
}

# Light-gray shading highlight on a few of the new "NewMotor" cells
foreach ($r in @(3, 6, 8, 9)) {
    $ws.Cells.Item($r, 5).Interior.Color = 15921906
}

# Column widths (best-fit like sizing for the new columns)
$ws.Columns.Item(5).ColumnWidth = 9.833333333333334
$ws.Columns.Item(6).ColumnWidth = 6.5

# --- View changes -------------------------------------------------------------
$ws.Range("F11").Select() | Out-Null
$excel.ActiveWindow.Zoom = 190
